$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("E2").Value = 25.89000000000061
$ws.Range("G2").Value = [double]"1.124504489524725e-09"
$ws.Range("H2").Value = [double]"4.029123701343451e-08"
$ws.Range("I2").Value = ""
$ws.Range("K2").Value = 5.979361791161358
$ws.Range("L2").Value = "[3.965814563305347, 7.992909019017368]"
$ws.Range("M2").Value = [double]"9.938334377679325e-09"
$ws.Range("N2").Value = [double]"1.987666875535865e-08"
$ws.Range("O2").Value = -1.509473947670772
$ws.Range("P2").Value = "[-1.9120003337163114, -1.1069475616252324]"
$ws.Range("Q2").Value = [double]"7.65609797781508e-13"
$ws.Range("R2").Value = [double]"1.531219595563016e-12"
$ws.Range("S2").Value = 11.55894877224376
$ws.Range("T2").Value = "[10.291853369044826, 12.8260441754427]"
$ws.Range("W2").Value = 6.219819819819968
$ws.Range("X2").Value = 4.561201201201309
$ws.Range("Y2").Value = 7.878438438438627

# --- Row 3 updates ---
$ws.Range("E3").Value = 22.28000000000004
$ws.Range("G3").Value = 0.01175142647344285
$ws.Range("H3").Value = 0.0401060366353803
$ws.Range("K3").Value = 4.495017184662181
$ws.Range("L3").Value = "[1.0116040804693966, 7.978430288854966]"
$ws.Range("M3").Value = 0.01171211008495754
$ws.Range("N3").Value = 0.01171211008495754
$ws.Range("O3").Value = 1.641552918091964
$ws.Range("P3").Value = "[0.5597632555945777, 2.7233425805893505]"
$ws.Range("Q3").Value = 0.003128618108784176
$ws.Range("R3").Value = 0.003128618108784176
$ws.Range("S3").Value = 12.73133883152051
$ws.Range("T3").Value = "[10.670840068125047, 14.791837594915966]"
$ws.Range("W3").Value = 16.45909909909913
$ws.Range("X3").Value = 12.62310310310313
$ws.Range("Y3").Value = 20.29509509509514
